# Apply updated values to "Generator Data" and "Yearly Fuel Costs" sheets

$wb = $excel.ActiveWorkbook

# --- Sheet: Generator Data ---
$ws1 = $wb.Worksheets.Item("Generator Data")

$ws1.Range("B2").Value = 7956.41204699
$ws1.Range("C2").Value = 29973.2158777

$ws1.Range("B3").Value = 1591.282409398
$ws1.Range("C3").Value = 7493.303969425

$ws1.Range("B4").Value = 71.60770842290999
$ws1.Range("C4").Value = 337.198678624125

$ws1.Range("B5").Value = 591.0851266240001
$ws1.Range("C5").Value = 23359.4267404

# --- Sheet: Yearly Fuel Costs ---
$ws2 = $wb.Worksheets.Item("Yearly Fuel Costs")

$ws2.Range("B2").Value = 75.72664166297506
$ws2.Range("B3").Value = 960.6693361899218
$ws2.Range("B4").Value = 983.9071927690284
$ws2.Range("B5").Value = 13329.86551201737
$ws2.Range("B6").Value = 29948.96724055574
